$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 corresponds to report "GF1103-121-境内汇总数据-季-人民币"
# Update development status from "未开始" (Not started) to "完成" (Completed),
# matching the look already used for "完成" status cells elsewhere in column C
# (green "Good" text/fill, keeping the existing thin-border formatting).
$ws.Range("C3").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("C19").Value = "完成"
$excel.CutCopyMode = 0

# Move the active selection to D1 (as recorded in the saved view state)
$ws.Range("D1").Select()
